$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 currently shows the text label "R40". It needs to become the
# text "1" instead, while keeping its existing cell style (s="23") and
# being stored as a genuine text value (not an auto-converted number).
#
# A direct `$ws.Range("B11").Value = "1"` would be auto-coerced by Excel
# into the *number* 1 (since "1" looks numeric and the cell's number
# format is General). To keep it a text value without touching the
# cell's number format/style, build the text in a scratch cell (using a
# text formula so it is unambiguously a string), copy it, and paste only
# the resulting value into B11 - this preserves B11's existing style.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163) # xlPasteValues
$scratch.Clear()
